$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(316, 'Marcos', 3, 1, 5, 1, 1, 0, 0, 0, 0, 0),
    @(317, 'Marcelão', 3, 1, 5, 0, 1, 0, 0, 0, 0, 0),
    @(318, 'Boneco', 3, 1, 5, 0, 1, 0, 0, 0, 0, 0),
    @(319, 'Romario', 3, 1, 5, 3, 1, 0, 0, 0, 0, 0),
    @(320, 'Coxinha', 3, 1, 5, 3, 1, 0, 0, 0, 0, 0),
    @(321, 'Eduardo', 4, 2, 3, 3, 1, 0, 0, 0, 0, 0),
    @(322, 'Du', 4, 2, 3, 1, 1, 0, 0, 0, 0, 0),
    @(323, 'Euler', 4, 2, 3, 0, 1, 0, 0, 0, 0, 0),
    @(324, 'Fernando', 4, 2, 3, 2, 1, 0, 0, 0, 0, 0),
    @(325, 'Leandrinho', 4, 2, 3, 1, 1, 0, 0, 0, 0, 0),
    @(326, 'Caio', 6, 1, 3, 4, 1, 1, 0, 0, 0, 0),
    @(327, 'Vander', 6, 1, 3, 1, 1, 1, 0, 0, 0, 0),
    @(328, 'Juscielio', 6, 1, 3, 3, 1, 1, 0, 1, 0, 0),
    @(329, 'Guinha', 6, 1, 3, 0, 1, 1, 0, 0, 0, 0),
    @(330, 'Michel', 6, 1, 3, 2, 1, 1, 0, 0, 0, 0),
    @(331, 'Corinthiano', 3, 0, 5, 0, 1, 0, 1, 0, 0, 0),
    @(332, 'Cabeleira', 3, 0, 5, 1, 1, 0, 1, 0, 0, 0),
    @(333, 'Athos', 3, 0, 5, 2, 1, 0, 1, 0, 0, 0),
    @(334, 'Eder', 3, 0, 5, 1, 1, 0, 1, 0, 0, 0),
    @(335, 'Alan', 3, 0, 5, 2, 1, 0, 1, 0, 0, 0),
    @(336, 'Matheus', 6, 1, 5, 0, 1, 0, 0, 0, 8, 1),
    @(337, 'Chelin', 7, 2, 4, 1, 1, 1, 0, 0, 8, 0),
    @(338, 'Lucian', 2, 1, 6, 0, 1, 0, 1, 0, 11, 0),
)

foreach ($row in $rows) {
    $r = $row[0]
    $name = $row[1]
    $ws.Cells.Item($r, 1).Value = $name
    for ($c = 0; $c -lt 10; $c++) {
        $ws.Cells.Item($r, 3 + $c).Value = $row[2 + $c]
    }
}

$ws.Range("A338").Select() | Out-Null
